# Generate Report for Handback
#
# The 63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.md file has now been handed back
# (status moves from "In Translation" to "Handed back: in sync with en-US",
# with fresh handback timestamps and the stale "error detail" note cleared).
# The report re-sorts so the row that was just refreshed lands on top of
# each table, pushing the other two (otherwise-unchanged) rows down by one.
#
# NOTE: literal "True"/"False" text values are written with a leading
# apostrophe so Excel stores them as text (matching the source data) instead
# of auto-coercing them into native Boolean cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.md"
$ov.Range("B2").Value = "e2e\63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.md"
$ov.Range("C2").Value = ".md"
$ov.Range("D2").Value = ""
$ov.Range("E2").Value = "Handed back: in sync with en-US"
$ov.Range("F2").Value = "Handed back: in sync with en-US"
$ov.Range("G2").Value = "2016-10-21 04:24:40"

$ov.Range("A3").Value = "ffff3f2694a0-01a4-436c-9d39-a13344643e0b.md"
$ov.Range("B3").Value = "e2e\ffff3f2694a0-01a4-436c-9d39-a13344643e0b.md"
$ov.Range("C3").Value = ".md"
$ov.Range("D3").Value = ""
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"
$ov.Range("G3").Value = "2016-10-21 04:19:19"

$ov.Range("A4").Value = "ffffff5cc3e3b1-7d27-4d27-855e-d51f152b3724.md"
$ov.Range("B4").Value = "e2e\ffffff5cc3e3b1-7d27-4d27-855e-d51f152b3724.md"
$ov.Range("C4").Value = ".md"
$ov.Range("D4").Value = ""
$ov.Range("E4").Value = "Handed back: in sync with en-US"
$ov.Range("F4").Value = "Handed back: in sync with en-US"
$ov.Range("G4").Value = "2016-10-21 04:19:19"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc256c0b853190a2bf123e915a6885c70f4fa264/e2e/ffff3f2694a0-01a4-436c-9d39-a13344643e0b.md", "", "", "e2e\63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.md")
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/32a4ccaa915593c3e2253a61f7a2e7492a0c3b01/e2e/ffffff5cc3e3b1-7d27-4d27-855e-d51f152b3724.md", "", "", "e2e\ffff3f2694a0-01a4-436c-9d39-a13344643e0b.md")
$ov.Hyperlinks.Add($ov.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/32a4ccaa915593c3e2253a61f7a2e7492a0c3b01/e2e/63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.md", "", "", "e2e\ffffff5cc3e3b1-7d27-4d27-855e-d51f152b3724.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("D2").Value = "e2e"
$zh.Range("E2").Value = "ht"
$zh.Range("F2").Value = "'False"
$zh.Range("G2").Value = "63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.1f7c8658640fd5f1cbde591f3e59de9d76f3ab1e.zh-cn.xlf"
$zh.Range("H2").Value = "2016-10-21 04:24:28"
$zh.Range("I2").Value = "63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.md"
$zh.Range("J2").Value = "63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.1f7c8658640fd5f1cbde591f3e59de9d76f3ab1e.zh-cn.xlf"
$zh.Range("K2").Value = "2016-10-21 04:25:07"
$zh.Range("L2").Value = ""
$zh.Range("M2").Value = "'True"
$zh.Range("N2").Value = ""
$zh.Range("O2").Value = "'False"
$zh.Range("P2").Value = ""

$zh.Range("A3").Value = "ffff3f2694a0-01a4-436c-9d39-a13344643e0b.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("D3").Value = "e2e"
$zh.Range("E3").Value = "ht"
$zh.Range("F3").Value = "'False"
$zh.Range("G3").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.c2607544a66b02a746f17728b9b1fcccf78d1073.zh-cn.xlf"
$zh.Range("H3").Value = "2016-10-21 04:19:08"
$zh.Range("I3").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.md"
$zh.Range("J3").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.c2607544a66b02a746f17728b9b1fcccf78d1073.zh-cn.xlf"
$zh.Range("K3").Value = "2016-10-21 04:19:48"
$zh.Range("L3").Value = ""
$zh.Range("M3").Value = "'True"
$zh.Range("N3").Value = ""
$zh.Range("O3").Value = "'False"
$zh.Range("P3").Value = ""

$zh.Range("A4").Value = "ffffff5cc3e3b1-7d27-4d27-855e-d51f152b3724.md"
$zh.Range("B4").Value = ".md"
$zh.Range("C4").Value = "Handed back: in sync with en-US"
$zh.Range("D4").Value = "e2e"
$zh.Range("E4").Value = "ht"
$zh.Range("F4").Value = "'True"
$zh.Range("G4").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.c2607544a66b02a746f17728b9b1fcccf78d1073.zh-cn.xlf"
$zh.Range("H4").Value = "2016-10-21 04:19:08"
$zh.Range("I4").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.md"
$zh.Range("J4").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.c2607544a66b02a746f17728b9b1fcccf78d1073.zh-cn.xlf"
$zh.Range("K4").Value = "2016-10-21 04:19:48"
$zh.Range("L4").Value = ""
$zh.Range("M4").Value = "'True"
$zh.Range("N4").Value = ""
$zh.Range("O4").Value = "'False"
$zh.Range("P4").Value = ""

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc256c0b853190a2bf123e915a6885c70f4fa264/e2e/ffff3f2694a0-01a4-436c-9d39-a13344643e0b.md", "", "", "63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.md")
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9c0aacd6f8ad7f84e745c1057021f2ae03c96cc7/e2e/09572edd-dbb6-4c5e-ac11-fa5758def696.md", "", "", "63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/32a4ccaa915593c3e2253a61f7a2e7492a0c3b01/e2e/ffffff5cc3e3b1-7d27-4d27-855e-d51f152b3724.md", "", "", "ffff3f2694a0-01a4-436c-9d39-a13344643e0b.md")
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9c0aacd6f8ad7f84e745c1057021f2ae03c96cc7/e2e/09572edd-dbb6-4c5e-ac11-fa5758def696.md", "", "", "09572edd-dbb6-4c5e-ac11-fa5758def696.md")
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/32a4ccaa915593c3e2253a61f7a2e7492a0c3b01/e2e/63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.md", "", "", "ffffff5cc3e3b1-7d27-4d27-855e-d51f152b3724.md")
$zh.Hyperlinks.Add($zh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/82d42b6b140b48ad6860ae33ec06f4150119b3a3/e2e/63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.md", "", "", "09572edd-dbb6-4c5e-ac11-fa5758def696.md")

$zh.Columns.Item(16).ColumnWidth = 13

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("D2").Value = "e2e"
$de.Range("E2").Value = "ht"
$de.Range("F2").Value = "'False"
$de.Range("G2").Value = "63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.1f7c8658640fd5f1cbde591f3e59de9d76f3ab1e.de-de.xlf"
$de.Range("H2").Value = "2016-10-21 04:24:40"
$de.Range("I2").Value = "63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.md"
$de.Range("J2").Value = "63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.1f7c8658640fd5f1cbde591f3e59de9d76f3ab1e.de-de.xlf"
$de.Range("K2").Value = "2016-10-21 04:25:25"
$de.Range("L2").Value = ""
$de.Range("M2").Value = "'True"
$de.Range("N2").Value = ""
$de.Range("O2").Value = "'False"
$de.Range("P2").Value = ""

$de.Range("A3").Value = "ffff3f2694a0-01a4-436c-9d39-a13344643e0b.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("D3").Value = "e2e"
$de.Range("E3").Value = "ht"
$de.Range("F3").Value = "'False"
$de.Range("G3").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.c2607544a66b02a746f17728b9b1fcccf78d1073.de-de.xlf"
$de.Range("H3").Value = "2016-10-21 04:19:19"
$de.Range("I3").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.md"
$de.Range("J3").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.c2607544a66b02a746f17728b9b1fcccf78d1073.de-de.xlf"
$de.Range("K3").Value = "2016-10-21 04:20:07"
$de.Range("L3").Value = ""
$de.Range("M3").Value = "'True"
$de.Range("N3").Value = ""
$de.Range("O3").Value = "'False"
$de.Range("P3").Value = ""

$de.Range("A4").Value = "ffffff5cc3e3b1-7d27-4d27-855e-d51f152b3724.md"
$de.Range("B4").Value = ".md"
$de.Range("C4").Value = "Handed back: in sync with en-US"
$de.Range("D4").Value = "e2e"
$de.Range("E4").Value = "ht"
$de.Range("F4").Value = "'True"
$de.Range("G4").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.c2607544a66b02a746f17728b9b1fcccf78d1073.de-de.xlf"
$de.Range("H4").Value = "2016-10-21 04:19:19"
$de.Range("I4").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.md"
$de.Range("J4").Value = "09572edd-dbb6-4c5e-ac11-fa5758def696.c2607544a66b02a746f17728b9b1fcccf78d1073.de-de.xlf"
$de.Range("K4").Value = "2016-10-21 04:20:07"
$de.Range("L4").Value = ""
$de.Range("M4").Value = "'True"
$de.Range("N4").Value = ""
$de.Range("O4").Value = "'False"
$de.Range("P4").Value = ""

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc256c0b853190a2bf123e915a6885c70f4fa264/e2e/ffff3f2694a0-01a4-436c-9d39-a13344643e0b.md", "", "", "63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.md")
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f469a79c4792e06fd404142a79a60c3b8b307f05/e2e/09572edd-dbb6-4c5e-ac11-fa5758def696.md", "", "", "63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.md")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/32a4ccaa915593c3e2253a61f7a2e7492a0c3b01/e2e/ffffff5cc3e3b1-7d27-4d27-855e-d51f152b3724.md", "", "", "ffff3f2694a0-01a4-436c-9d39-a13344643e0b.md")
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f469a79c4792e06fd404142a79a60c3b8b307f05/e2e/09572edd-dbb6-4c5e-ac11-fa5758def696.md", "", "", "09572edd-dbb6-4c5e-ac11-fa5758def696.md")
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/32a4ccaa915593c3e2253a61f7a2e7492a0c3b01/e2e/63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.md", "", "", "ffffff5cc3e3b1-7d27-4d27-855e-d51f152b3724.md")
$de.Hyperlinks.Add($de.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c780f92ddb72d4b7ef2fd81c484eabb417e110ba/e2e/63d6e77c-51c8-42be-942e-6bbb7eaa3bf3.md", "", "", "09572edd-dbb6-4c5e-ac11-fa5758def696.md")

$de.Columns.Item(16).ColumnWidth = 13
